$d = $word.ActiveDocument

$pairs = @(
    @("441×4=", "126×7="),
    @("380×2=", "467×3="),
    @("818×5=", "806×2="),
    @("300×6=", "147×9="),
    @("984×9=", "346×3="),
    @("442×7=", "756×7="),
    @("472×5=", "516×3="),
    @("217×5=", "158×3="),
    @("879×5=", "779×9="),
    @("969×5=", "369×3="),
    @("349×7=", "992×9="),
    @("842×5=", "788×9="),
    @("176×7=", "526×6="),
    @("846×2=", "878×9="),
    @("823×7=", "847×3="),
    @("509×6=", "782×5="),
    @("546×6=", "850×6="),
    @("179×9=", "223×2="),
    @("275×5=", "252×7="),
    @("780×8=", "485×8="),
    @("470×8=", "218×9="),
    @("811×7=", "823×9="),
    @("240×8=", "508×8="),
    @("383×5=", "848×2="),
    @("947×3=", "844×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
